$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (ColumnWidth adds ~0.8333 padding when stored in OOXML,
# so subtract 5/6 to land on the target stored width)
$ws.Columns.Item(6).ColumnWidth = 12.1666666666667
$ws.Columns.Item(7).ColumnWidth = 12.1666666666667
$ws.Columns.Item(8).ColumnWidth = 14.1666666666667

# Update header cell values
$ws.Range("F1").Value = "input_phone"
$ws.Range("G1").Value = "input_state"
$ws.Range("H1").Value = "input_zipCode"
